$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.250.44"
$ws.Range("E2").Value = "  +0.93%  "

$ws.Range("D3").Value = "1.564.61"
$ws.Range("E3").Value = "  +0.20%  "

$ws.Range("E4").Value = "  -0.30%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.33%  "

$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("E7").Value = "  -0.24%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.13"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("E9").Value = "  -0.18%  "

$ws.Range("E11").Value = "  +1.90%  "

$ws.Range("D12").Value = "1.787.49"
$ws.Range("E12").Value = "  +0.22%  "

$ws.Range("D13").Value = "1.565.58"
$ws.Range("E13").Value = "  +1.26%  "

$ws.Range("E14").Value = "  -0.09%  "

$ws.Range("D16").Value = "27.205.35"
$ws.Range("E16").Value = "  +0.78%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.08%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.77%  "

$ws.Range("E19").Value = "  +1.08%  "

$ws.Range("E20").Value = "  -0.80%  "

$ws.Range("E21").Value = "  -0.22%  "

$ws.Range("E22").Value = "  +0.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.38%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.29%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.72%  "

$ws.Range("E26").Value = "  +0.07%  "

$ws.Range("E27").Value = "  +1.24%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.66%  "

$ws.Range("E29").Value = "  -0.38%  "

$ws.Range("E30").Value = "  +1.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0470"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.90%  "

$ws.Range("E32").Value = "  +0.08%  "

$ws.Range("D33").Value = "1.455.90"

$ws.Range("E34").Value = "  +0.55%  "

$ws.Range("E35").Value = "  +4.61%  "

$ws.Range("E36").Value = "  +1.25%  "

$ws.Range("E37").Value = "  +0.73%  "

$ws.Range("E38").Value = "  -0.19%  "

$ws.Range("E39").Value = "  +0.88%  "

$ws.Range("E40").Value = "  +0.36%  "

$ws.Range("E41").Value = "  +0.62%  "

$ws.Range("E42").Value = "  -0.22%  "

$ws.Range("E43").Value = "  +1.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.981"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.72%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.69%  "

$ws.Range("E46").Value = "  +0.21%  "

$ws.Range("D47").Value = "1.699.39"
$ws.Range("E47").Value = "  -0.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.78%  "

$ws.Range("E49").Value = "  -1.62%  "

$ws.Range("E50").Value = "  +1.32%  "

$ws.Range("E51").Value = "  -1.79%  "
